$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.81 = 6666.67 pesos`n✅ 6666.67 pesos = 1.8 = 934.66 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 553.5
$wsTasas.Range("O10").Value = 3690
$wsTasas.Range("N12").Value = 3709
$wsTasas.Range("O12").Value = 520
